# Apply the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.422.72"
$ws.Range("E2").Value = "  +0.63%  "

$ws.Range("D3").Value = "3.975.52"
$ws.Range("E3").Value = "  -0.99%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.94"
$ws.Range("E5").Value = "  +10.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.90"
$ws.Range("E6").Value = "  +0.93%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.679"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.749"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("E10").Value = "  -2.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.74"
$ws.Range("E11").Value = "  +5.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000318"
$ws.Range("E12").Value = "  -2.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.84"
$ws.Range("E13").Value = "  +1.33%  "

$ws.Range("D14").Value = "4.629.20"
$ws.Range("E14").Value = "  -0.62%  "

$ws.Range("D15").Value = "3.984.45"
$ws.Range("E15").Value = "  -1.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.26"
$ws.Range("E16").Value = "  +6.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.98"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.46"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("E19").Value = "  -0.63%  "

$ws.Range("D20").Value = "72.406.92"
$ws.Range("E20").Value = "  +0.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "433.71"
$ws.Range("E21").Value = "  +0.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.73"
$ws.Range("E22").Value = "  +13.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "95.95"
$ws.Range("E23").Value = "  -1.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.49"
$ws.Range("E24").Value = "  -0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.25"
$ws.Range("E25").Value = "  -0.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.40"
$ws.Range("E26").Value = "  +19.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.57"
$ws.Range("E27").Value = "  +4.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.63"
$ws.Range("E28").Value = "  -0.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.91"
$ws.Range("E29").Value = "  +1.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.27"
$ws.Range("E30").Value = "  -0.93%  "

$ws.Range("E31").Value = "  +7.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.56"
$ws.Range("E32").Value = "  +4.99%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.132"
$ws.Range("E33").Value = "  +0.91%  "

$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.52"
$ws.Range("E34").Value = "  +0.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "681.05"
$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "69.49"
$ws.Range("E36").Value = "  +6.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.440"
$ws.Range("E37").Value = "  -1.09%  "

$ws.Range("D38").Value = "0.0₃0851"
$ws.Range("E38").Value = "  +3.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.147"
$ws.Range("E39").Value = "  -2.66%  "

$ws.Range("E40").Value = "  +0.31%  "

$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.07"
$ws.Range("E42").Value = "  +7.97%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.18%  "

$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.30"
$ws.Range("E44").Value = "  -2.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0486"
$ws.Range("E45").Value = "  -0.57%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.81"
$ws.Range("E46").Value = "  +5.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.149"
$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.36"
$ws.Range("E48").Value = "  -1.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.46"
$ws.Range("E49").Value = "  +6.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.99"
$ws.Range("E50").Value = "  -0.40%  "

$ws.Range("D51").Value = "2.778.57"
$ws.Range("E51").Value = "  +8.95%  "
